# issue #5: stock data output to json file
#
# Adds a "property_category" column (value "stock") to the 股票 (stock)
# worksheet, between the existing "total" and "date" columns, shifting the
# later columns (date, legislator_name, legislator_id) one column to the
# right (H -> I, I -> J, J -> K) and filling the vacated H column with the
# new property_category/"stock" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# xlPasteValues = -4163 (copies only the literal value/type, no re-parsing,
# so date-looking text like "2011-11-21" is not reinterpreted as a date
# serial number) and xlPasteFormats = -4122 (copies only cell formatting).
$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Header row (row 1): shift K1<-J1, J1<-I1, I1<-H1, then set H1 ---
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteValues)

$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial($xlPasteValues)

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial($xlPasteValues)

$ws.Range("H1").Value = "property_category"

# --- Data row (row 2): shift K2<-J2, J2<-I2, I2<-H2, then set H2 ---
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial($xlPasteValues)

$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial($xlPasteValues)

$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial($xlPasteValues)

$ws.Range("H2").Value = "stock"

# --- Fix up formatting on the newly-introduced column K so it matches the
#     rest of the table (K1 = bold/centered/bordered header style like the
#     other header cells, K2 = plain data style like the other data cells).
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteFormats)

$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
